# "Added in empty soft body scene"
#
# - TARGET PRACTISE status (H7) moves from UNDERWAY to a brand new "SATISFACTORY"
#   status: light green fill (FF92D050) with a thin box border (like the other
#   category-header status cells).
# - "Make Target Scene" (H8) moves from UNDERWAY to DONE, reusing the existing
#   DONE look (green fill FF00B050, no border).
# - A new "COMPOUND SHAPES" category/header row is appended (row 17), styled
#   like the other category headers, with a TODO status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 17: COMPOUND SHAPES category header, empty notes, TODO status
# (written first so the new shared strings land in the same order as the
# author's edit: "COMPOUND SHAPES" before "SATISFACTORY")
$ws.Range("F15:H15").Copy()
$ws.Range("F17:H17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F17").Value = "COMPOUND SHAPES"
$ws.Range("H17").Value = "TODO"

# --- H7: TARGET PRACTISE status -> SATISFACTORY (new light-green + border style)
$h7 = $ws.Range("H7")
$h7.Interior.Color = 5296274   # BGR for RGB FF92D050
$h7.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$h7.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$h7.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$h7.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$h7.Value = "SATISFACTORY"

# --- H8: Make Target Scene status -> DONE (reuse existing green/no-border look)
$h8 = $ws.Range("H8")
$h8.Interior.Color = 5287936   # BGR for RGB FF00B050
$h8.Value = "DONE"

# --- New column width tweak for column H (slightly wider to fit "SATISFACTORY")
$ws.Range("H1").EntireColumn.ColumnWidth = 15

# --- Update the saved selection/active cell to match the author's last position
$ws.Range("I24").Select()
